$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values
$ws.Range("A1").Value = "참가부"
$ws.Range("B1").Value = "체급"
$ws.Range("C1").Value = "인원수"

# Clear old column D and E header
$ws.Range("D1").Value = ""
$ws.Range("E1").Value = ""

# Remove the now-obsolete data rows (rows 2 and 3)
$ws.Range("A2:E3").Clear()
